# Table_S5.xlsx edit: "typo in read counts"
# - Renames D1 header "Final abundances3" -> "Mapped reads3"
# - Adds new E column "Final abundances4" with mapped/normalized read counts
# - Adds a "Mean:" row (row 7) with AVERAGE formulas across B:E
# - Rewrites the footnote block (now 5 entries, renumbered + reworded),
#   moved down to rows 10-14
# - Fixes the superscript footnote ref on the 408920 figure (4 -> 5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Header row: fix D1 text, add E1
# ---------------------------------------------------------------
$ws.Range("D1").Value = "Mapped reads3"
$d1base = $ws.Range("D1").Characters(1, 12)
$d1base.Font.Bold = $true
$d1base.Font.Name = "Arial"
$d1base.Font.Size = 10
$d1sup = $ws.Range("D1").Characters(13, 1)
$d1sup.Font.Bold = $true
$d1sup.Font.Superscript = $true
$d1sup.Font.Name = "Arial"
$d1sup.Font.Size = 10

$ws.Range("E1").Value = "Final abundances4"
$ws.Range("E1").Borders.LineStyle = 1
$e1base = $ws.Range("E1").Characters(1, 16)
$e1base.Font.Bold = $true
$e1base.Font.Name = "Arial"
$e1base.Font.Size = 10
$e1sup = $ws.Range("E1").Characters(17, 1)
$e1sup.Font.Bold = $true
$e1sup.Font.Superscript = $true
$e1sup.Font.Name = "Arial"
$e1sup.Font.Size = 10

# ---------------------------------------------------------------
# 2. New column E data rows 2-5 ("Final abundances" figures)
# ---------------------------------------------------------------
$ws.Range("E2").Value = 299325
$ws.Range("E3").Value = 222811
$ws.Range("E4").Value = 108783

$ws.Range("E2:E4").Borders.LineStyle = 1
$ws.Range("E2:E4").Font.Name = "Arial"
$ws.Range("E2:E4").Font.Size = 10

$ws.Range("E5").Value = 37829
$ws.Range("E5").Borders.LineStyle = 1
$ws.Range("E5").Font.Name = "Arial"
$ws.Range("E5").Font.Size = 10
$ws.Range("E5").HorizontalAlignment = -4152   # xlRight

# ---------------------------------------------------------------
# 3. Fix the footnote superscript on the D5 figure: 408920^4 -> 408920^5
#    (force text storage first so the numeric-looking string isn't
#    silently coerced to a Double before we can split it into runs)
# ---------------------------------------------------------------
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "4089205"
$ws.Range("D5").NumberFormat = "General"
$d5base = $ws.Range("D5").Characters(1, 6)
$d5base.Font.Name = "Arial"
$d5base.Font.Size = 10
$d5sup = $ws.Range("D5").Characters(7, 1)
$d5sup.Font.Bold = $true
$d5sup.Font.Superscript = $true
$d5sup.Font.Name = "Arial"
$d5sup.Font.Size = 10

# ---------------------------------------------------------------
# 4. Clear the old footnote block (rows 7-10, incl. their merges)
#    and rebuild: row 7 becomes "Mean:" with AVERAGE formulas, and
#    the renumbered / reworded footnotes move down to rows 10-14.
# ---------------------------------------------------------------
$ws.Range("A7:I14").UnMerge()
$ws.Range("A7:I14").Clear()

$ws.Range("A7").Value = "Mean:"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Name = "Arial"
$ws.Range("A7").Font.Size = 10
$ws.Range("A7").Borders.LineStyle = 1

$ws.Range("B7").Formula = "=AVERAGE(B2:B5)"
$ws.Range("C7").Formula = "=AVERAGE(C2:C5)"
$ws.Range("D7").Formula = "=AVERAGE(D2:D5)"
$ws.Range("E7").Formula = "=AVERAGE(E2:E5)"
$ws.Range("B7:E7").NumberFormat = "#"
$ws.Range("B7:E7").Font.Bold = $true
$ws.Range("B7:E7").Font.Name = "Arial"
$ws.Range("B7:E7").Font.Size = 10
$ws.Range("B7:E7").Borders.LineStyle = 1

# Footnote 1
$ws.Range("A10").Value = "1 = Raw paired-end reads from sequencer, includes all metatranscriptomic reads"
$ws.Range("A10").Font.Name = "Arial"
$ws.Range("A10").Font.Size = 10
$ws.Range("A10:C10").Merge()

# Footnote 2
$ws.Range("A11").Value = "2 = Totals include both paired-end and orphaned reads following quality and adapter trimming, includes all metatranscriptomic reads"
$ws.Range("A11").Font.Name = "Arial"
$ws.Range("A11").Font.Size = 10
$ws.Range("A11:E11").Merge()

# Footnote 3 (rich text: italic species name)
$ws.Range("A12").Value = "3 = Mapped reads to C. difficile strain 630 genome from Bowtie2 from each pool"
$a12pre = $ws.Range("A12").Characters(1, 20)
$a12pre.Font.Name = "Arial"
$a12pre.Font.Size = 10
$a12italic = $ws.Range("A12").Characters(21, 12)
$a12italic.Font.Italic = $true
$a12italic.Font.Name = "Arial"
$a12italic.Font.Size = 10
$a12post = $ws.Range("A12").Characters(33, 46)
$a12post.Font.Name = "Arial"
$a12post.Font.Size = 10
$ws.Range("A12:C12").Merge()

# Footnote 4 (new)
$ws.Range("A13").Value = "4 = Read counts after removal of optical + PCR duplicates and normalization to target gene/read length"
$ws.Range("A13").Font.Name = "Arial"
$ws.Range("A13").Font.Size = 10
$ws.Range("A13:D13").Merge()

# Footnote 5 (formerly footnote 4, reworded)
$ws.Range("A14").Value = "5 = Library with RIN score 6.8, likely caused fewer reads mapped to reference genome with no mismatches"
$ws.Range("A14").Font.Name = "Arial"
$ws.Range("A14").Font.Size = 10
$ws.Range("A14:D14").Merge()

# ---------------------------------------------------------------
# 5. Misc cosmetic bits captured in the diff
# ---------------------------------------------------------------
$wb.Windows.Item(1).TabRatio = 0.99
$ws.Range("D20").Select()
